$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow the "Recorded By" column (G) now that its contents are being redacted.
$ws.Columns.Item(7).ColumnWidth = 12.17

# Redact the "Recorded By" name for every session that was actually recorded.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $status = $ws.Cells.Item($r, 9).Value2
    if ($status -eq "Recorded") {
        $ws.Cells.Item($r, 7).Value = ""
    }
}
